$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @{ B = 0.02258322285507441; C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 2.328612170846122 }
    3 = @{ B = 1.505614041169197;   C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 3.811642989160245 }
    4 = @{ B = 1.505614041169197;   C = 1.65323645889881;  D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
    5 = @{ B = 1.505614041169197;   C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 3.811642989160245 }
    6 = @{ B = 3.182878228561681;   C = 1.65323645889881;  D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    7 = @{ B = 0.7287194209349384;  C = 0.3375848360084654; D = 3.082599426703578;  E = 0.4998867070740569; G = 4.64879039072104 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
